# fix for busybox run cmd in notes
#
# Slide 9's speaker notes describe the "kubectl run" demo command. The
# command was missing the new pod/container name ("dns-test") and still
# referenced the (deprecated in newer busybox images) "/bin/sh" shell
# instead of "/bin/ash". Update the notes text accordingly:
#   "kubectl run --rm -ti --image=busybox /bin/sh"
# becomes
#   "kubectl run dns-test --rm -ti --image=busybox /bin/ash"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$notesShape = $s.NotesPage.Shapes.Item(2)
$tr = $notesShape.TextFrame.TextRange

$tr.Text = "Demo: get a pod + shell session with " + [char]8220 + "kubectl run dns-test --rm -ti --image=busybox /bin/ash" + [char]8221 + "; use the DNS name of a service to download an index.html (i.e. " + [char]8220 + "wget nginx" + [char]8221 + ")"
